$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("R1").Value = "Cooldown"

$values = @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,2,8,1,2,0,3,2,2,0,0)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 18).Value = $values[$i]
}
